# Updates power-flow result values (pl_mw) for the 380 kV case (Case_5_62).
# Columns B,C,D,E,F,H,J,K,M,N are updated for rows 2-25 (index 0-23);
# columns A,G,I,L,O are left untouched (G/I/L/O stay 0, A is the row index).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("B", "C", "D", "E", "F", "H", "J", "K", "M", "N")

$newValues = @{
    2 = @(1.040936523214157, 0.1440678763281937, 0.03147556612577418, 0.09419238487479475, 3.472502459719323, 0.07973214163530429, 0.1843251799547474, 1.106647970836462, 0.3884562439538328, 3.608084878746723)
    3 = @(1.002616785891206, 0.1380182215032733, 0.03117916304550405, 0.09372680780215248, 3.45337420963719, 0.07973214163530429, 0.1836150428678209, 1.064745862659947, 0.3789740781784801, 3.612708044780433)
    4 = @(0.9797222471347311, 0.134397175520732, 0.03101622485735689, 0.09348286438032893, 3.44299368482848, 0.07973214163530429, 0.1832525573497392, 1.039701203945185, 0.3733789150369162, 3.61631952472267)
    5 = @(0.9705518261576458, 0.1329450154809848, 0.03095464164427497, 0.09339401245362744, 3.439106426931801, 0.07973214163530429, 0.1831233378920061, 1.029666963550596, 0.3711558974215379, 3.617985388974688)
    6 = @(0.9690387021591391, 0.1327052998250053, 0.03094470734930255, 0.09337989672776992, 3.438481655583487, 0.07973214163530429, 0.1831029983534407, 1.028011148267694, 0.370790213034752, 3.618273725886809)
    7 = @(0.9795979268280064, 0.1343774963281845, 0.0310153747947588, 0.09348162332666377, 3.442939871802082, 0.07973214163530429, 0.1832507397546053, 1.039565183862692, 0.3733487036287855, 3.616341205275361)
    8 = @(1.027592216632485, 0.141962519283851, 0.03136942572535872, 0.09402316149296297, 3.465623759463284, 0.07973214163530429, 0.1840650626497649, 1.092058152974346, 0.385139687066939, 3.609518432470026)
    9 = @(1.126750821028764, 0.1575819356264105, 0.03221394087285034, 0.09541725817975433, 3.520945794967957, 0.07973214163530429, 0.1862455951658433, 1.200434798318099, 0.4100644082094647, 3.60228062256445)
    10 = @(1.202702450768129, 0.1695182526035808, 0.03292481664951197, 0.09664356862873547, 3.568226683539621, 0.07973214163530429, 0.1882040475268383, 1.28340624885476, 0.4294811042292395, 3.600722327107007)
    11 = @(1.23793421016336, 0.1750499290531025, 0.03326760813297369, 0.0972452591954216, 3.591183626566931, 0.07973214163530429, 0.189172556746577, 1.321886528387523, 0.4385554914805283, 3.600833217581851)
    12 = @(1.251373807403183, 0.1771593605086821, 0.03340018388646371, 0.09747939734703337, 3.600085530014979, 0.07973214163530429, 0.1895504710763802, 1.33656428018611, 0.4420265348665779, 3.600993353464872)
    13 = @(1.24847498410395, 0.1767044015521151, 0.03337150858218507, 0.09742869191006776, 3.59815906314185, 0.07973214163530429, 0.1894685840394317, 1.333398440637723, 0.4412774367781509, 3.600953606608172)
    14 = @(1.239037926955518, 0.1752231781253784, 0.03327845987536193, 0.09726439587205249, 3.59191180898992, 0.07973214163530429, 0.1892034242903762, 1.323091947911365, 0.4388403593288501, 3.600844022754117)
    15 = @(1.233270236942928, 0.174317803185204, 0.03322182465624479, 0.09716457871037321, 3.588112361220396, 0.07973214163530429, 0.1890424598412892, 1.316792749018617, 0.4373521078082234, 3.600792293197671)
    16 = @(1.200413689924517, 0.1691587985933438, 0.03290280302495319, 0.09660512774504681, 3.566755568557227, 0.07973214163530429, 0.1881423148476102, 1.280906306296345, 0.4288929359154068, 3.600731598292796)
    17 = @(1.180431797296251, 0.1660200348091223, 0.0327120485374266, 0.09627314108555041, 3.55402510454347, 0.07973214163530429, 0.1876099815405468, 1.259079855439296, 0.423765412163533, 3.600904511160891)
    18 = @(1.169002827821032, 0.1642242772678912, 0.03260415998094146, 0.09608631790544564, 3.546839211980014, 0.07973214163530429, 0.1873111013667312, 1.246595132885801, 0.4208389384247155, 3.601081110432887)
    19 = @(1.165144177003441, 0.1636179072941388, 0.03256794550429021, 0.09602377194925182, 3.544429598773178, 0.07973214163530429, 0.1872111600611746, 1.242379910926076, 0.4198519900443429, 3.601154145594677)
    20 = @(1.182552269948019, 0.1663531699589385, 0.03273216561541403, 0.09630805461937442, 3.555366172084732, 0.07973214163530429, 0.187665893393536, 1.261396146965325, 0.424308891919722, 3.600878118300599)
    21 = @(1.241807155747722, 0.1756578496876102, 0.03330571559492768, 0.09731248296392181, 3.593741115018133, 0.07973214163530429, 0.1892810052124361, 1.326116332984753, 0.4395552442122366, 3.600873001674017)
    22 = @(1.281105503327353, 0.1818247445084467, 0.03369668745625631, 0.09800559804657993, 3.620037426138339, 0.07973214163530429, 0.1904016298295375, 1.36903329808689, 0.4497222945052712, 3.601558417073221)
    23 = @(1.260078852868105, 0.1785254875142357, 0.03348655067276951, 0.09763231887786716, 3.605891219436586, 0.07973214163530429, 0.1897975779276706, 1.346071021915321, 0.4442773952003805, 3.601129488439227)
    24 = @(1.181593420852039, 0.1662025322797263, 0.03272306514195122, 0.09629225762489568, 3.554759461120042, 0.07973214163530429, 0.187640593326357, 1.260348753498477, 0.4240631180417154, 3.600889810081881)
    25 = @(1.099383343701561, 0.153276030053803, 0.03196951659367997, 0.09500460391917187, 3.504816464229151, 0.07973214163530429, 0.1855931587324875, 1.170530588199938, 0.4031280583613253, 3.603579582798488)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $vals[$i]
    }
}

Write-Output "Updated $($newValues.Keys.Count) rows across $($columns.Count) columns"
